$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "27.564.75"
$ws.Cells.Item(2, 5).Value = "  -0.17%  "

$ws.Cells.Item(3, 4).Value = "1.835.70"
$ws.Cells.Item(3, 5).Value = "  -0.47%  "

$ws.Cells.Item(4, 4).Value = "'1.002"
$ws.Cells.Item(4, 5).Value = "  -0.07%  "

$ws.Cells.Item(5, 4).Value = "'312.27"
$ws.Cells.Item(5, 5).Value = "  -0.18%  "

$ws.Cells.Item(6, 5).Value = "  -0.10%  "

$ws.Cells.Item(7, 4).Value = "'0.4278"
$ws.Cells.Item(7, 5).Value = "  -0.14%  "

$ws.Cells.Item(8, 4).Value = "'0.3653"
$ws.Cells.Item(8, 5).Value = "  +0.49%  "

$ws.Cells.Item(9, 4).Value = "'0.07269"
$ws.Cells.Item(9, 5).Value = "  -0.66%  "

$ws.Cells.Item(10, 4).Value = "'0.8646"
$ws.Cells.Item(10, 5).Value = "  -1.79%  "

$ws.Cells.Item(11, 4).Value = "'20.67"
$ws.Cells.Item(11, 5).Value = "  +0.07%  "

$ws.Cells.Item(12, 2).Value = "Polkadot"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Cells.Item(12, 4).Value = "'5.469"
$ws.Cells.Item(12, 5).Value = "  +2.29%  "

$ws.Cells.Item(13, 2).Value = "WrappedEther"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(13, 4).Value = "1.740.13"
$ws.Cells.Item(13, 5).Value = "  -9.27%  "

$ws.Cells.Item(14, 4).Value = "'6.523"
$ws.Cells.Item(14, 5).Value = "  -0.15%  "

$ws.Cells.Item(15, 4).Value = "'0.06974"
$ws.Cells.Item(15, 5).Value = "  +0.32%  "

$ws.Cells.Item(16, 4).Value = "'1.002"
$ws.Cells.Item(16, 5).Value = "  -0.20%  "

$ws.Cells.Item(17, 4).Value = "'80.67"
$ws.Cells.Item(17, 5).Value = "  +1.37%  "

$ws.Cells.Item(18, 4).Value = "'0.000008891"
$ws.Cells.Item(18, 5).Value = "  -0.93%  "

$ws.Cells.Item(19, 4).Value = "'1.000"
$ws.Cells.Item(19, 5).Value = "  -0.11%  "

$ws.Cells.Item(20, 4).Value = "'15.42"
$ws.Cells.Item(20, 5).Value = "  +0.30%  "

$ws.Cells.Item(21, 4).Value = "27.523.25"
$ws.Cells.Item(21, 5).Value = "  -0.09%  "

$ws.Cells.Item(22, 4).Value = "'5.157"
$ws.Cells.Item(22, 5).Value = "  +3.28%  "

$ws.Cells.Item(23, 5).Value = "  +5.51%  "

$ws.Cells.Item(24, 4).Value = "2.036.31"
$ws.Cells.Item(24, 5).Value = "  -3.16%  "

$ws.Cells.Item(25, 4).Value = "'1.990"
$ws.Cells.Item(25, 5).Value = "  +0.12%  "

$ws.Cells.Item(26, 4).Value = "'154.87"
$ws.Cells.Item(26, 5).Value = "  -0.37%  "

$ws.Cells.Item(27, 4).Value = "'18.88"
$ws.Cells.Item(27, 5).Value = "  +1.89%  "

$ws.Cells.Item(28, 4).Value = "'5.165"
$ws.Cells.Item(28, 5).Value = "  -1.06%  "

$ws.Cells.Item(29, 4).Value = "'114.32"
$ws.Cells.Item(29, 5).Value = "  -4.58%  "

$ws.Cells.Item(30, 4).Value = "'1.820"
$ws.Cells.Item(30, 5).Value = "  -3.57%  "

$ws.Cells.Item(31, 4).Value = "'0.08855"
$ws.Cells.Item(31, 5).Value = "  -0.50%  "

$ws.Cells.Item(32, 4).Value = "'0.7498"
$ws.Cells.Item(32, 5).Value = "  -2.06%  "

$ws.Cells.Item(33, 4).Value = "'2.997"
$ws.Cells.Item(33, 5).Value = "  +1.05%  "

$ws.Cells.Item(34, 4).Value = "'4.550"
$ws.Cells.Item(34, 5).Value = "  +0.52%  "

$ws.Cells.Item(35, 4).Value = "'1.134"
$ws.Cells.Item(35, 5).Value = "  +0.26%  "

$ws.Cells.Item(36, 5).Value = "  -0.06%  "

$ws.Cells.Item(37, 4).Value = "'1.098"

$ws.Cells.Item(38, 4).Value = "'0.05324"
$ws.Cells.Item(38, 5).Value = "  -2.91%  "

$ws.Cells.Item(39, 4).Value = "'0.01938"
$ws.Cells.Item(39, 5).Value = "  -0.05%  "

$ws.Cells.Item(40, 4).Value = "'2.799"
$ws.Cells.Item(40, 5).Value = "  -0.71%  "

$ws.Cells.Item(41, 4).Value = "'0.5076"
$ws.Cells.Item(41, 5).Value = "  -0.06%  "

$ws.Cells.Item(42, 4).Value = "'0.1650"
$ws.Cells.Item(42, 5).Value = "  -1.04%  "

$ws.Cells.Item(43, 4).Value = "'6.473"
$ws.Cells.Item(43, 5).Value = "  -1.86%  "

$ws.Cells.Item(44, 4).Value = "'8.338"
$ws.Cells.Item(44, 5).Value = "  -0.96%  "

$ws.Cells.Item(45, 4).Value = "'10.46"
$ws.Cells.Item(45, 5).Value = "  +0.87%  "

$ws.Cells.Item(46, 4).Value = "'105.59"
$ws.Cells.Item(46, 5).Value = "  -0.29%  "

$ws.Cells.Item(47, 4).Value = "'0.06478"
$ws.Cells.Item(47, 5).Value = "  -1.05%  "

$ws.Cells.Item(48, 4).Value = "'0.4694"
$ws.Cells.Item(48, 5).Value = "  +0.84%  "

$ws.Cells.Item(49, 4).Value = "'1.0000"
$ws.Cells.Item(49, 5).Value = "  -0.14%  "

$ws.Cells.Item(50, 4).Value = "'1.621"
$ws.Cells.Item(50, 5).Value = "  -1.14%  "

$ws.Cells.Item(51, 4).Value = "'1.741"
$ws.Cells.Item(51, 5).Value = "  -0.27%  "
